# Auto-generated edit script applying numeric corrections to H:N columns
# across sheets ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR, per the target diff.
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 1781.421
$ws.Range("I19").Value = 3084.0557
$ws.Range("J19").Value = 609.05
$ws.Range("K19").Value = 3084.0557
$ws.Range("L19").Value = 609.05
$ws.Range("M19").Value = -2909.0557
$ws.Range("N19").Value = -959.05
$ws.Range("H92").Value = 3357.45
$ws.Range("I92").Value = 3991
$ws.Range("J92").Value = 1879.1666
$ws.Range("K92").Value = 3991
$ws.Range("L92").Value = 1879.1666
$ws.Range("M92").Value = -2743
$ws.Range("N92").Value = -4375.1666
$ws.Range("H96").Value = 1246.375
$ws.Range("I96").Value = 1149.3846
$ws.Range("J96").Value = 1666.6666
$ws.Range("K96").Value = 3448.1538
$ws.Range("L96").Value = 4999.9998
$ws.Range("M96").Value = -2075.1538
$ws.Range("N96").Value = -7745.9998
$ws.Range("H117").Value = 43000
$ws.Range("J117").Value = 43000
$ws.Range("L117").Value = 43000
$ws.Range("N117").Value = -52178
$ws.Range("H132").Value = 2724.2173
$ws.Range("I132").Value = 1711.6818
$ws.Range("K132").Value = 5135.0454
$ws.Range("M132").Value = -2605.0454
$ws.Range("H137").Value = 3206.3215
$ws.Range("I137").Value = 3224.0417
$ws.Range("J137").Value = 3100
$ws.Range("K137").Value = 9672.125100000001
$ws.Range("L137").Value = 9300
$ws.Range("M137").Value = -7122.125100000001
$ws.Range("N137").Value = -14400

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1019302.75
$ws.Range("I32").Value = 1228691.9
$ws.Range("J32").Value = 2269.5715
$ws.Range("K32").Value = 1228691.9
$ws.Range("L32").Value = 2269.5715
$ws.Range("M32").Value = -1228404.9
$ws.Range("N32").Value = -2843.5715
$ws.Range("H45").Value = 1240.6154
$ws.Range("I45").Value = 1102.5454
$ws.Range("K45").Value = 1102.5454
$ws.Range("M45").Value = -725.5454
$ws.Range("H61").Value = 574684.1
$ws.Range("I61").Value = 528612.75
$ws.Range("J61").Value = 629393.9
$ws.Range("K61").Value = 528612.75
$ws.Range("L61").Value = 629393.9
$ws.Range("M61").Value = -528400.75
$ws.Range("N61").Value = -629817.9
$ws.Range("H110").Value = 1704.5
$ws.Range("I110").Value = 1834.4445
$ws.Range("J110").Value = 1314.6666
$ws.Range("K110").Value = 1834.4445
$ws.Range("L110").Value = 1314.6666
$ws.Range("M110").Value = 210.5554999999999
$ws.Range("N110").Value = -5404.6666
$ws.Range("H111").Value = 16500
$ws.Range("J111").Value = 16500
$ws.Range("L111").Value = 16500
$ws.Range("N111").Value = -24680
$ws.Range("H122").Value = 1532.6333
$ws.Range("I122").Value = 1281.2727
$ws.Range("J122").Value = 2223.875
$ws.Range("K122").Value = 3843.8181
$ws.Range("L122").Value = 6671.625
$ws.Range("M122").Value = -1393.8181
$ws.Range("N122").Value = -11571.625
$ws.Range("H132").Value = 17015.637
$ws.Range("I132").Value = 22925.541
$ws.Range("K132").Value = 68776.62300000001
$ws.Range("M132").Value = -66246.62300000001
$ws.Range("H136").Value = 574684.1
$ws.Range("I136").Value = 528612.75
$ws.Range("J136").Value = 629393.9
$ws.Range("K136").Value = 1585838.25
$ws.Range("L136").Value = 1888181.7
$ws.Range("M136").Value = -1583288.25
$ws.Range("N136").Value = -1893281.7

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 927.8570999999999
$ws.Range("I20").Value = 854.5
$ws.Range("J20").Value = 1025.6666
$ws.Range("K20").Value = 854.5
$ws.Range("L20").Value = 1025.6666
$ws.Range("M20").Value = -607.5
$ws.Range("N20").Value = -1519.6666
$ws.Range("H107").Value = 1425.7727
$ws.Range("I107").Value = 1303.4117
$ws.Range("J107").Value = 1841.8
$ws.Range("K107").Value = 1303.4117
$ws.Range("L107").Value = 1841.8
$ws.Range("M107").Value = 616.5882999999999
$ws.Range("N107").Value = -5681.8
$ws.Range("H129").Value = 39332.668
$ws.Range("J129").Value = 39332.668
$ws.Range("L129").Value = 39332.668
$ws.Range("N129").Value = -49332.668
$ws.Range("H134").Value = 2099.8845
$ws.Range("I134").Value = 1163
$ws.Range("J134").Value = 4642.857
$ws.Range("K134").Value = 3489
$ws.Range("L134").Value = 13928.571
$ws.Range("M134").Value = -954
$ws.Range("N134").Value = -18998.571

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2941.8096
$ws.Range("I31").Value = 1968.875
$ws.Range("J31").Value = 6055.2
$ws.Range("K31").Value = 1968.875
$ws.Range("L31").Value = 6055.2
$ws.Range("M31").Value = -1673.875
$ws.Range("N31").Value = -6645.2
$ws.Range("H34").Value = 2941.8096
$ws.Range("I34").Value = 1968.875
$ws.Range("J34").Value = 6055.2
$ws.Range("K34").Value = 1968.875
$ws.Range("L34").Value = 6055.2
$ws.Range("M34").Value = -1766.875
$ws.Range("N34").Value = -6459.2
$ws.Range("H52").Value = 27780
$ws.Range("J52").Value = 27780
$ws.Range("L52").Value = 27780
$ws.Range("N52").Value = -28368
$ws.Range("H58").Value = 6458.909
$ws.Range("I58").Value = 10949.6
$ws.Range("K58").Value = 10949.6
$ws.Range("M58").Value = -10746.6
$ws.Range("H99").Value = 73993.07000000001
$ws.Range("I99").Value = 145287.86
$ws.Range("J99").Value = 2698.2856
$ws.Range("K99").Value = 145287.86
$ws.Range("L99").Value = 2698.2856
$ws.Range("M99").Value = -143789.86
$ws.Range("N99").Value = -5694.2856
$ws.Range("H122").Value = 3013.9
$ws.Range("I122").Value = 3243.2222
$ws.Range("J122").Value = 950
$ws.Range("K122").Value = 9729.6666
$ws.Range("L122").Value = 2850
$ws.Range("M122").Value = -7279.6666
$ws.Range("N122").Value = -7750
$ws.Range("H126").Value = 73993.07000000001
$ws.Range("I126").Value = 145287.86
$ws.Range("J126").Value = 2698.2856
$ws.Range("K126").Value = 435863.58
$ws.Range("L126").Value = 8094.8568
$ws.Range("M126").Value = -433393.58
$ws.Range("N126").Value = -13034.8568
$ws.Range("H132").Value = 1665.1282
$ws.Range("I132").Value = 1018.06665
$ws.Range("J132").Value = 3822
$ws.Range("K132").Value = 3054.19995
$ws.Range("L132").Value = 11466
$ws.Range("M132").Value = -524.1999500000002
$ws.Range("N132").Value = -16526
$ws.Range("H134").Value = 1283.8286
$ws.Range("I134").Value = 785.9259
$ws.Range("J134").Value = 2964.25
$ws.Range("K134").Value = 2357.7777
$ws.Range("L134").Value = 8892.75
$ws.Range("M134").Value = 177.2223000000004
$ws.Range("N134").Value = -13962.75
$ws.Range("H136").Value = 6458.909
$ws.Range("I136").Value = 10949.6
$ws.Range("K136").Value = 32848.8
$ws.Range("M136").Value = -30298.8

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 1082.4
$ws.Range("I14").Value = 1082.4
$ws.Range("K14").Value = 3247.2
$ws.Range("M14").Value = -3074.2
$ws.Range("H20").Value = 2999.6667
$ws.Range("J20").Value = 2999.6667
$ws.Range("L20").Value = 8999.000100000001
$ws.Range("N20").Value = -9453.000100000001
$ws.Range("H82").Value = 7750
$ws.Range("J82").Value = 7750
$ws.Range("L82").Value = 23250
$ws.Range("N82").Value = -24062
$ws.Range("H85").Value = 7750
$ws.Range("J85").Value = 7750
$ws.Range("L85").Value = 23250
$ws.Range("N85").Value = -26058
$ws.Range("H126").Value = 7379.3
$ws.Range("I126").Value = 8837.143
$ws.Range("J126").Value = 3977.6667
$ws.Range("K126").Value = 26511.429
$ws.Range("L126").Value = 11933.0001
$ws.Range("M126").Value = -21571.429
$ws.Range("N126").Value = -21813.0001
$ws.Range("H132").Value = 4510.737
$ws.Range("I132").Value = 1492
$ws.Range("J132").Value = 9685.714
$ws.Range("K132").Value = 13428
$ws.Range("L132").Value = 87171.42600000001
$ws.Range("M132").Value = -10898
$ws.Range("N132").Value = -92231.42600000001

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 369.64285
$ws.Range("I107").Value = 355.22223
$ws.Range("K107").Value = 355.22223
$ws.Range("M107").Value = 1564.77777
$ws.Range("H122").Value = 1439.5
$ws.Range("I122").Value = 1439.5
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 4318.5
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -1868.5
$ws.Range("N122").ClearContents()
$ws.Range("H132").Value = 4921.3667
$ws.Range("I132").Value = 6909.364
$ws.Range("K132").Value = 20728.092
$ws.Range("M132").Value = -18198.092

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 2900
$ws.Range("I40").Value = 2875
$ws.Range("J40").Value = 3000
$ws.Range("K40").Value = 2875
$ws.Range("L40").Value = 3000
$ws.Range("M40").Value = -2739
$ws.Range("N40").Value = -3272
$ws.Range("H61").Value = 1306.4615
$ws.Range("I61").Value = 1312
$ws.Range("J61").Value = 1300
$ws.Range("K61").Value = 1312
$ws.Range("L61").Value = 1300
$ws.Range("M61").Value = -1110
$ws.Range("N61").Value = -1704
$ws.Range("H113").Value = 1306.4615
$ws.Range("I113").Value = 1312
$ws.Range("J113").Value = 1300
$ws.Range("K113").Value = 1312
$ws.Range("L113").Value = 1300
$ws.Range("M113").Value = 858
$ws.Range("N113").Value = -5640
$ws.Range("H129").Value = 35666.668
$ws.Range("I129").Value = 0
$ws.Range("J129").Value = 35666.668
$ws.Range("K129").Value = 0
$ws.Range("L129").Value = 35666.668
$ws.Range("M129").ClearContents()
$ws.Range("N129").Value = -45666.668
$ws.Range("H136").Value = 4471
$ws.Range("I136").Value = 3013.7896
$ws.Range("J136").Value = 13700
$ws.Range("K136").Value = 9041.3688
$ws.Range("L136").Value = 41100
$ws.Range("M136").Value = -6491.3688
$ws.Range("N136").Value = -46200

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 71014.5
$ws.Range("J46").Value = 71014.5
$ws.Range("L46").Value = 71014.5
$ws.Range("N46").Value = -71476.5
$ws.Range("H107").Value = 475.75
$ws.Range("I107").Value = 398.52942
$ws.Range("K107").Value = 1195.58826
$ws.Range("M107").Value = 724.41174
$ws.Range("H126").Value = 1491.5834
$ws.Range("I126").Value = 1187.375
$ws.Range("J126").Value = 2100
$ws.Range("K126").Value = 3562.125
$ws.Range("L126").Value = 6300
$ws.Range("M126").Value = -1092.125
$ws.Range("N126").Value = -11240
$ws.Range("H129").Value = 36750
$ws.Range("J129").Value = 36750
$ws.Range("L129").Value = 36750
$ws.Range("N129").Value = -46750
$ws.Range("H134").Value = 71014.5
$ws.Range("J134").Value = 71014.5
$ws.Range("L134").Value = 213043.5
$ws.Range("N134").Value = -218113.5
$ws.Range("H136").Value = 18383736
$ws.Range("I136").Value = 24416156
$ws.Range("J136").Value = 717363.9399999999
$ws.Range("K136").Value = 73248468
$ws.Range("L136").Value = 2152091.82
$ws.Range("M136").Value = -73245918
$ws.Range("N136").Value = -2157191.82

